$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-8 from serial date 45175 to 45183
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45183
}
